$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Checklist "X" marks moved between columns ---
# Row 29: mark moves from E (N/A) to C (OK)
$ws.Range("C29").Value = "X"
$ws.Range("E29").Value = ""

# Row 33: mark removed from C (OK)
$ws.Range("C33").Value = ""

# Row 42: mark removed from C (OK)
$ws.Range("C42").Value = ""

# Row 58: mark added to C (OK)
$ws.Range("C58").Value = "X"

# Row 60: mark added to E (N/A)
$ws.Range("E60").Value = "X"

# --- Comments column (F) updates ---
$ws.Range("F30").Value = "Still needs final revision, some may not be acording to standard"
$ws.Range("F32").Value = "?"
$ws.Range("F33").Value = "Power width of 0.508mm and all other tracks with 0.254mm?"
$ws.Range("F40").Value = "Not a high speed design"
$ws.Range("F41").Value = "Simple design that doesn't require it"
$ws.Range("F42").Value = "?"
$ws.Range("F43").Value = "?"
$ws.Range("F45").Value = "?"
$ws.Range("F48").Value = "Initial placement but is not showing hole in drill drawing layer"
$ws.Range("F47").Value = "too large table and it is not know how to resize it"
$ws.Range("F54").Value = "?"
$ws.Range("F55").Value = "?"
$ws.Range("F56").Value = "?"
$ws.Range("F57").Value = "Simple design that doesn't require it"
$ws.Range("F59").Value = "?"
$ws.Range("F60").Value = "Simple design that doesn't require it"
$ws.Range("F61").Value = "?"

# --- Column F width adjustment (auto-fit after longer comments were entered) ---
$ws.Columns("F:F").ColumnWidth = 59.46

# --- View state: scroll position & active selection ---
$excel.ActiveWindow.ScrollRow = 43
$ws.Range("D56").Select()
